# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型"
# sheets, refreshed by the gh-pages data regeneration (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibition) sheet ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 433   # 合肥·W·A第五人格同人only2.0            432 -> 433
$wsExpo.Range("F3").Value = 7     # 安徽·崩坏同人only 爱莉希雅同人生日会    5 -> 7
$wsExpo.Range("F4").Value = 3291  # 合肥·第九届环形宇宙动漫游戏嘉年华      3265 -> 3291
$wsExpo.Range("F5").Value = 154   # 合肥·MAX特摄同人only2.0               153 -> 154
$wsExpo.Range("F7").Value = 152   # 合肥·心动恋章·冬日序国乙&代号鸢同人only 144 -> 152

# --- 全部类型 (All types) sheet ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 433    # 合肥·W·A第五人格同人only2.0            432 -> 433
$wsAll.Range("F7").Value = 7      # 安徽·崩坏同人only 爱莉希雅同人生日会    5 -> 7
$wsAll.Range("F8").Value = 3291   # 合肥·第九届环形宇宙动漫游戏嘉年华      3265 -> 3291
$wsAll.Range("F9").Value = 154    # 合肥·MAX特摄同人only2.0               153 -> 154
$wsAll.Range("F12").Value = 152   # 合肥·心动恋章·冬日序国乙&代号鸢同人only 144 -> 152
